$d = $word.ActiveDocument
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(3)
$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Step 3: diabetes prediction through rule based</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:t>import pandas as pd</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">import </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>plotly.express</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> as </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>px</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
        <w:t># Step 1: Load dataset</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">data = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>pd.read_csv</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>("diabetes.csv")</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>print(" Data loaded successfully!")</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>print(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>data.head</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>())</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
        <w:t># Step 2: Define simple rule-based function</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">def </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>predict_diabetes</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(row):</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    # Simple logic — not ML, just conditions</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    if (</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">        row["Glucose"] &gt; 130</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">        or row["</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>BloodPressure</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>"] &gt; 85</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">        or row["BMI"] &gt; 30</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">        or row["Age"] &gt; 45</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    ):</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">        return 1  # Diabetic</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    else:</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">        return 0  # Non-diabetic</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
        <w:t># Step 3: Apply logic to each row</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>data["</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Predicted_Outcome</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">"] = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>data.apply</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>predict_diabetes</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, axis=1)</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
        <w:t># Step 4: Evaluate</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>correct = (data["</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Predicted_Outcome</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>"] == data["Outcome"]).sum()</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">accuracy = correct / </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>len</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(data) * 100</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>print(f"\n Simple Rule-Based Accuracy: {accuracy:.2f}%")</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:br/>
        <w:t># Step 5: Create bar chart for Actual vs Predicted</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>summary = (</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    data[["Outcome", "</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Predicted_Outcome</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>"]]</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    .melt(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>var_name</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">="Type", </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>value_name</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>="Value")</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    .</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>groupby</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(["Type", "Value"])</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    .size()</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    .</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>reset_index</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(name="Count")</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
        <w:t># Step 6: Visualization</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">fig = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>px.bar</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    summary,</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    x="Value",</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    y="Count",</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    color="Type",</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>barmode</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>="group",</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    text="Count",</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>color_discrete_sequence</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>=["#1f77b4", "#ff7f0e"],</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    title="Actual vs Predicted Diabetes Cases (Rule-Based)"</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>fig.update_traces</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>texttemplate</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">='%{text}', </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>textposition</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>='outside')</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>fig.update_layout</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>xaxis_title</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>="Diabetes (0 = No, 1 = Yes)",</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>yaxis_title</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>="Number of Patients",</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>title_font</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>=</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>dict</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(size=20, color="</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>darkblue</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>"),</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">    template="</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>plotly_white</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>fig.show</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>()</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>print("\n Visualization displayed successfully!")</w:t>
      </w:r>
    </w:p>

'@
$newPara.Range.InsertXML($xml)
